$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(2, 8).Value = 1230.8
$ws.Cells.Item(2, 9).Value = 718
$ws.Cells.Item(2, 11).Value = 718
$ws.Cells.Item(2, 13).Value = -605

$ws.Cells.Item(9, 8).Value = 429038.44
$ws.Cells.Item(9, 9).Value = 600349.9
$ws.Cells.Item(9, 10).Value = 759.75
$ws.Cells.Item(9, 11).Value = 600349.9
$ws.Cells.Item(9, 12).Value = 759.75
$ws.Cells.Item(9, 13).Value = -600180.9
$ws.Cells.Item(9, 14).Value = -1097.75

$ws.Cells.Item(15, 8).Value = 1196.7297
$ws.Cells.Item(15, 9).Value = 1196.7297
$ws.Cells.Item(15, 11).Value = 3590.189100000001
$ws.Cells.Item(15, 13).Value = -3421.189100000001

$ws.Cells.Item(17, 8).Value = 373254.28
$ws.Cells.Item(17, 10).Value = 373254.28
$ws.Cells.Item(17, 12).Value = 1119762.84
$ws.Cells.Item(17, 14).Value = -1120098.84

$ws.Cells.Item(28, 8).Value = 980.2692
$ws.Cells.Item(28, 10).Value = 1959.8
$ws.Cells.Item(28, 12).Value = 1959.8
$ws.Cells.Item(28, 14).Value = -2929.8

$ws.Cells.Item(34, 8).Value = 7374.5713
$ws.Cells.Item(34, 9).Value = 7374.5713
$ws.Cells.Item(34, 11).Value = 7374.5713
$ws.Cells.Item(34, 13).Value = -7171.5713

$ws.Cells.Item(36, 8).Value = 7374.5713
$ws.Cells.Item(36, 9).Value = 7374.5713
$ws.Cells.Item(36, 11).Value = 7374.5713
$ws.Cells.Item(36, 13).Value = -6659.5713

$ws.Cells.Item(40, 8).Value = 4624.6665
$ws.Cells.Item(40, 10).Value = 4199
$ws.Cells.Item(40, 12).Value = 4199
$ws.Cells.Item(40, 14).Value = -4549

$ws.Cells.Item(53, 8).Value = 1262.92
$ws.Cells.Item(53, 9).Value = 91.111115
$ws.Cells.Item(53, 10).Value = 1922.0625
$ws.Cells.Item(53, 11).Value = 91.111115
$ws.Cells.Item(53, 12).Value = 1922.0625
$ws.Cells.Item(53, 13).Value = 545.888885
$ws.Cells.Item(53, 14).Value = -3196.0625

$ws.Cells.Item(62, 8).Value = 2128.4092
$ws.Cells.Item(62, 9).Value = 2383.2
$ws.Cells.Item(62, 10).Value = 1916.0834
$ws.Cells.Item(62, 11).Value = 2383.2
$ws.Cells.Item(62, 12).Value = 1916.0834
$ws.Cells.Item(62, 13).Value = -1759.2
$ws.Cells.Item(62, 14).Value = -3164.0834

$ws.Cells.Item(64, 8).Value = 5401
$ws.Cells.Item(64, 10).Value = 8000
$ws.Cells.Item(64, 12).Value = 8000
$ws.Cells.Item(64, 14).Value = -8496

$ws.Cells.Item(65, 8).Value = 2128.4092
$ws.Cells.Item(65, 9).Value = 2383.2
$ws.Cells.Item(65, 10).Value = 1916.0834
$ws.Cells.Item(65, 11).Value = 11916
$ws.Cells.Item(65, 12).Value = 9580.416999999999
$ws.Cells.Item(65, 13).Value = -8796
$ws.Cells.Item(65, 14).Value = -15820.417

$ws.Cells.Item(67, 8).Value = 5401
$ws.Cells.Item(67, 10).Value = 8000
$ws.Cells.Item(67, 12).Value = 8000
$ws.Cells.Item(67, 14).Value = -9716

$ws.Cells.Item(70, 8).Value = 54781.81
$ws.Cells.Item(70, 9).Value = 101730
$ws.Cells.Item(70, 10).Value = 12101.637
$ws.Cells.Item(70, 11).Value = 305190
$ws.Cells.Item(70, 12).Value = 36304.911
$ws.Cells.Item(70, 13).Value = -304920
$ws.Cells.Item(70, 14).Value = -36844.911

$ws.Cells.Item(73, 8).Value = 54781.81
$ws.Cells.Item(73, 9).Value = 101730
$ws.Cells.Item(73, 10).Value = 12101.637
$ws.Cells.Item(73, 11).Value = 305190
$ws.Cells.Item(73, 12).Value = 36304.911
$ws.Cells.Item(73, 13).Value = -304254
$ws.Cells.Item(73, 14).Value = -38176.911

$ws.Cells.Item(80, 8).Value = 2449.3333
$ws.Cells.Item(80, 9).Value = 628.8333
$ws.Cells.Item(80, 10).Value = 3359.5833
$ws.Cells.Item(80, 11).Value = 1886.4999
$ws.Cells.Item(80, 12).Value = 10078.7499
$ws.Cells.Item(80, 13).Value = -888.4999
$ws.Cells.Item(80, 14).Value = -12074.7499

$ws.Cells.Item(83, 8).Value = 2449.3333
$ws.Cells.Item(83, 9).Value = 628.8333
$ws.Cells.Item(83, 10).Value = 3359.5833
$ws.Cells.Item(83, 11).Value = 5659.4997
$ws.Cells.Item(83, 12).Value = 30236.2497
$ws.Cells.Item(83, 13).Value = -667.4997000000003
$ws.Cells.Item(83, 14).Value = -40220.2497

$ws.Cells.Item(86, 8).Value = 2395
$ws.Cells.Item(86, 9).Value = 2239.2
$ws.Cells.Item(86, 11).Value = 2239.2
$ws.Cells.Item(86, 13).Value = -1116.2

$ws.Cells.Item(89, 8).Value = 2395
$ws.Cells.Item(89, 9).Value = 2239.2
$ws.Cells.Item(89, 11).Value = 11196
$ws.Cells.Item(89, 13).Value = -5580

$ws.Cells.Item(92, 8).Value = 1517.1428
$ws.Cells.Item(92, 9).Value = 1290.9166
$ws.Cells.Item(92, 11).Value = 1290.9166
$ws.Cells.Item(92, 13).Value = -42.91660000000002

$ws.Cells.Item(96, 8).Value = 657.6
$ws.Cells.Item(96, 9).Value = 767.6
$ws.Cells.Item(96, 10).Value = 437.6
$ws.Cells.Item(96, 11).Value = 2302.8
$ws.Cells.Item(96, 12).Value = 1312.8
$ws.Cells.Item(96, 13).Value = -929.8000000000002
$ws.Cells.Item(96, 14).Value = -4058.8

$ws.Cells.Item(98, 8).Value = 4752.9375
$ws.Cells.Item(98, 9).Value = 2045
$ws.Cells.Item(98, 11).Value = 2045
$ws.Cells.Item(98, 13).Value = -547

$ws.Cells.Item(101, 8).Value = 646.4167
$ws.Cells.Item(101, 9).Value = 508.625
$ws.Cells.Item(101, 10).Value = 922
$ws.Cells.Item(101, 11).Value = 1525.875
$ws.Cells.Item(101, 12).Value = 2766
$ws.Cells.Item(101, 13).Value = 96.125
$ws.Cells.Item(101, 14).Value = -6010

$ws.Cells.Item(106, 8).Value = 3566.6667
$ws.Cells.Item(106, 9).Value = 2850
$ws.Cells.Item(106, 11).Value = 2850
$ws.Cells.Item(106, 13).Value = -2219

$ws.Cells.Item(107, 8).Value = 982.17145
$ws.Cells.Item(107, 9).Value = 611.04
$ws.Cells.Item(107, 11).Value = 611.04
$ws.Cells.Item(107, 13).Value = 1308.96

$ws.Cells.Item(113, 8).Value = 5442.4287
$ws.Cells.Item(113, 9).Value = 2716.8333
$ws.Cells.Item(113, 11).Value = 2716.8333
$ws.Cells.Item(113, 13).Value = 537.1667000000002

$ws.Cells.Item(116, 8).Value = 2771.4443
$ws.Cells.Item(116, 9).Value = 2580.25
$ws.Cells.Item(116, 10).Value = 2924.4
$ws.Cells.Item(116, 11).Value = 2580.25
$ws.Cells.Item(116, 12).Value = 2924.4
$ws.Cells.Item(116, 13).Value = 861.75
$ws.Cells.Item(116, 14).Value = -9808.4

$ws.Cells.Item(122, 8).Value = 4752.9375
$ws.Cells.Item(122, 9).Value = 2045
$ws.Cells.Item(122, 11).Value = 6135
$ws.Cells.Item(122, 13).Value = -3685

$ws.Cells.Item(125, 8).Value = 6297.4
$ws.Cells.Item(125, 9).Value = 3500
$ws.Cells.Item(125, 10).Value = 6608.222
$ws.Cells.Item(125, 11).Value = 31500
$ws.Cells.Item(125, 12).Value = 59473.998
$ws.Cells.Item(125, 14).Value = -64393.998
$ws.Cells.Item(125, 13).Value = -29040

$ws.Cells.Item(126, 8).Value = 83597.78
$ws.Cells.Item(126, 10).Value = 83597.78
$ws.Cells.Item(126, 12).Value = 83597.78
$ws.Cells.Item(126, 14).Value = -93477.78

$ws.Cells.Item(132, 8).Value = 10536220
$ws.Cells.Item(132, 9).Value = 11372247
$ws.Cells.Item(132, 10).Value = 26171
$ws.Cells.Item(132, 11).Value = 34116741
$ws.Cells.Item(132, 12).Value = 78513
$ws.Cells.Item(132, 13).Value = -34114211
$ws.Cells.Item(132, 14).Value = -83573

$ws.Cells.Item(135, 8).Value = 2266.5
$ws.Cells.Item(135, 9).Value = 1473.7174
$ws.Cells.Item(135, 10).Value = 5913.3
$ws.Cells.Item(135, 11).Value = 13263.4566
$ws.Cells.Item(135, 12).Value = 53219.7
$ws.Cells.Item(135, 13).Value = -10728.4566
$ws.Cells.Item(135, 14).Value = -58289.7

$ws.Cells.Item(137, 8).Value = 5632.091
$ws.Cells.Item(137, 9).Value = 9151.666999999999
$ws.Cells.Item(137, 10).Value = 4312.25
$ws.Cells.Item(137, 11).Value = 27455.001
$ws.Cells.Item(137, 12).Value = 12936.75
$ws.Cells.Item(137, 13).Value = -24905.001
$ws.Cells.Item(137, 14).Value = -18036.75

$ws.Cells.Item(138, 8).Value = 581518.4399999999
$ws.Cells.Item(138, 9).Value = 126310.625
$ws.Cells.Item(138, 11).Value = 378931.875
$ws.Cells.Item(138, 13).Value = -373791.875

$ws.Cells.Item(141, 8).Value = 1091.5111
$ws.Cells.Item(141, 9).Value = 641.35
$ws.Cells.Item(141, 11).Value = 1924.05
$ws.Cells.Item(141, 13).Value = 3255.95

$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(2, 8).Value = 731.6667
$ws.Cells.Item(2, 9).Value = 565.05884
$ws.Cells.Item(2, 10).Value = 1439.75
$ws.Cells.Item(2, 11).Value = 565.05884
$ws.Cells.Item(2, 12).Value = 1439.75
$ws.Cells.Item(2, 13).Value = -452.05884
$ws.Cells.Item(2, 14).Value = -1665.75

$ws.Cells.Item(28, 8).Value = 15106.2
$ws.Cells.Item(28, 10).Value = 55500
$ws.Cells.Item(28, 12).Value = 55500
$ws.Cells.Item(28, 14).Value = -55884

$ws.Cells.Item(31, 8).Value = 971.4286
$ws.Cells.Item(31, 9).Value = 971.4286
$ws.Cells.Item(31, 10).Value = 0
$ws.Cells.Item(31, 11).Value = 971.4286
$ws.Cells.Item(31, 12).Value = 0
$ws.Cells.Item(31, 13).Value = -677.4286
$ws.Cells.Item(31, 14).ClearContents()

$ws.Cells.Item(32, 8).Value = 2790.4177
$ws.Cells.Item(32, 9).Value = 2714.5417
$ws.Cells.Item(32, 11).Value = 2714.5417
$ws.Cells.Item(32, 13).Value = -2427.5417

$ws.Cells.Item(41, 8).Value = 1280.5
$ws.Cells.Item(41, 9).Value = 1280.5
$ws.Cells.Item(41, 11).Value = 1280.5
$ws.Cells.Item(41, 13).Value = -866.5

$ws.Cells.Item(43, 8).Value = 17652.834
$ws.Cells.Item(43, 9).Value = 12670.5
$ws.Cells.Item(43, 10).Value = 20144
$ws.Cells.Item(43, 11).Value = 12670.5
$ws.Cells.Item(43, 12).Value = 20144
$ws.Cells.Item(43, 13).Value = -12357.5
$ws.Cells.Item(43, 14).Value = -20770

$ws.Cells.Item(45, 8).Value = 2856.0557
$ws.Cells.Item(45, 9).Value = 2442.5
$ws.Cells.Item(45, 10).Value = 3683.1667
$ws.Cells.Item(45, 11).Value = 2442.5
$ws.Cells.Item(45, 12).Value = 3683.1667
$ws.Cells.Item(45, 13).Value = -2065.5
$ws.Cells.Item(45, 14).Value = -4437.1667

$ws.Cells.Item(61, 8).Value = 5581.7827
$ws.Cells.Item(61, 9).Value = 4694.05
$ws.Cells.Item(61, 11).Value = 4694.05
$ws.Cells.Item(61, 13).Value = -4482.05

$ws.Cells.Item(70, 8).Value = 72900
$ws.Cells.Item(70, 10).Value = 72900
$ws.Cells.Item(70, 12).Value = 72900
$ws.Cells.Item(70, 14).Value = -73440

$ws.Cells.Item(73, 8).Value = 72900
$ws.Cells.Item(73, 10).Value = 72900
$ws.Cells.Item(73, 12).Value = 72900
$ws.Cells.Item(73, 14).Value = -74772

$ws.Cells.Item(74, 8).Value = 1805.2572
$ws.Cells.Item(74, 9).Value = 1334.5
$ws.Cells.Item(74, 10).Value = 3688.2856
$ws.Cells.Item(74, 11).Value = 1334.5
$ws.Cells.Item(74, 12).Value = 3688.2856
$ws.Cells.Item(74, 13).Value = -460.5
$ws.Cells.Item(74, 14).Value = -5436.2856

$ws.Cells.Item(77, 8).Value = 1805.2572
$ws.Cells.Item(77, 9).Value = 1334.5
$ws.Cells.Item(77, 10).Value = 3688.2856
$ws.Cells.Item(77, 11).Value = 6672.5
$ws.Cells.Item(77, 12).Value = 18441.428
$ws.Cells.Item(77, 13).Value = -2304.5
$ws.Cells.Item(77, 14).Value = -27177.428

$ws.Cells.Item(88, 8).Value = 1027.6666
$ws.Cells.Item(88, 9).Value = 535
$ws.Cells.Item(88, 10).Value = 1421.8
$ws.Cells.Item(88, 11).Value = 535
$ws.Cells.Item(88, 12).Value = 1421.8
$ws.Cells.Item(88, 13).Value = -129
$ws.Cells.Item(88, 14).Value = -2233.8

$ws.Cells.Item(91, 8).Value = 1027.6666
$ws.Cells.Item(91, 9).Value = 535
$ws.Cells.Item(91, 10).Value = 1421.8
$ws.Cells.Item(91, 11).Value = 535
$ws.Cells.Item(91, 12).Value = 1421.8
$ws.Cells.Item(91, 13).Value = 869
$ws.Cells.Item(91, 14).Value = -4229.8

$ws.Cells.Item(93, 8).Value = 68888
$ws.Cells.Item(93, 10).Value = 68888
$ws.Cells.Item(93, 12).Value = 68888
$ws.Cells.Item(93, 14).Value = -73880

$ws.Cells.Item(99, 8).Value = 15106.2
$ws.Cells.Item(99, 10).Value = 55500
$ws.Cells.Item(99, 12).Value = 55500
$ws.Cells.Item(99, 14).Value = -61490

$ws.Cells.Item(102, 8).Value = 4257480
$ws.Cells.Item(102, 9).Value = 2234.0789
$ws.Cells.Item(102, 11).Value = 2234.0789
$ws.Cells.Item(102, 13).Value = -612.0789

$ws.Cells.Item(104, 8).Value = 0
$ws.Cells.Item(104, 10).Value = 0
$ws.Cells.Item(104, 12).Value = 0
$ws.Cells.Item(104, 14).ClearContents()

$ws.Cells.Item(106, 8).Value = 82862.336
$ws.Cells.Item(106, 10).Value = 82862.336
$ws.Cells.Item(106, 12).Value = 82862.336
$ws.Cells.Item(106, 14).Value = -85386.336

$ws.Cells.Item(110, 8).Value = 1396.1482
$ws.Cells.Item(110, 9).Value = 943.26086
$ws.Cells.Item(110, 10).Value = 4000.25
$ws.Cells.Item(110, 11).Value = 943.26086
$ws.Cells.Item(110, 12).Value = 4000.25
$ws.Cells.Item(110, 13).Value = 1101.73914
$ws.Cells.Item(110, 14).Value = -8090.25

$ws.Cells.Item(116, 8).Value = 731.6667
$ws.Cells.Item(116, 9).Value = 565.05884
$ws.Cells.Item(116, 10).Value = 1439.75
$ws.Cells.Item(116, 11).Value = 565.05884
$ws.Cells.Item(116, 12).Value = 1439.75
$ws.Cells.Item(116, 13).Value = 1728.94116
$ws.Cells.Item(116, 14).Value = -6027.75

$ws.Cells.Item(122, 8).Value = 1537.125
$ws.Cells.Item(122, 9).Value = 1537.125
$ws.Cells.Item(122, 11).Value = 4611.375
$ws.Cells.Item(122, 13).Value = -2161.375

$ws.Cells.Item(132, 8).Value = 4200.383
$ws.Cells.Item(132, 9).Value = 2297.513
$ws.Cells.Item(132, 11).Value = 6892.539
$ws.Cells.Item(132, 13).Value = -4362.539

$ws.Cells.Item(136, 8).Value = 5581.7827
$ws.Cells.Item(136, 9).Value = 4694.05
$ws.Cells.Item(136, 11).Value = 14082.15
$ws.Cells.Item(136, 13).Value = -11532.15

$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(3, 8).Value = 731.6667
$ws.Cells.Item(3, 9).Value = 565.05884
$ws.Cells.Item(3, 10).Value = 1439.75
$ws.Cells.Item(3, 11).Value = 565.05884
$ws.Cells.Item(3, 12).Value = 1439.75
$ws.Cells.Item(3, 13).Value = -451.05884
$ws.Cells.Item(3, 14).Value = -1667.75

$ws.Cells.Item(15, 8).Value = 9007
$ws.Cells.Item(15, 10).Value = 9007
$ws.Cells.Item(15, 12).Value = 9007
$ws.Cells.Item(15, 14).Value = -9461

$ws.Cells.Item(86, 8).Value = 4764157.5
$ws.Cells.Item(86, 9).Value = 6063201
$ws.Cells.Item(86, 10).Value = 999
$ws.Cells.Item(86, 11).Value = 6063201
$ws.Cells.Item(86, 12).Value = 999
$ws.Cells.Item(86, 13).Value = -6062078
$ws.Cells.Item(86, 14).Value = -3245

$ws.Cells.Item(89, 8).Value = 4764157.5
$ws.Cells.Item(89, 9).Value = 6063201
$ws.Cells.Item(89, 10).Value = 999
$ws.Cells.Item(89, 11).Value = 30316005
$ws.Cells.Item(89, 12).Value = 4995
$ws.Cells.Item(89, 13).Value = -30310389
$ws.Cells.Item(89, 14).Value = -16227

$ws.Cells.Item(99, 8).Value = 3101.7778
$ws.Cells.Item(99, 9).Value = 2309.5715
$ws.Cells.Item(99, 11).Value = 2309.5715
$ws.Cells.Item(99, 13).Value = -811.5715

$ws.Cells.Item(105, 8).Value = 3576.238
$ws.Cells.Item(105, 9).Value = 2900.0527
$ws.Cells.Item(105, 11).Value = 2900.0527
$ws.Cells.Item(105, 13).Value = -1153.0527

$ws.Cells.Item(134, 8).Value = 5945.7812
$ws.Cells.Item(134, 9).Value = 2034.4166
$ws.Cells.Item(134, 10).Value = 8292.6
$ws.Cells.Item(134, 11).Value = 6103.2498
$ws.Cells.Item(134, 12).Value = 24877.8
$ws.Cells.Item(134, 13).Value = -3568.2498
$ws.Cells.Item(134, 14).Value = -29947.8

$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(4, 8).Value = 1000000000
$ws.Cells.Item(4, 10).Value = 1000000000
$ws.Cells.Item(4, 12).Value = 1000000000
$ws.Cells.Item(4, 14).Value = -1000000224

$ws.Cells.Item(22, 8).Value = 997.7143
$ws.Cells.Item(22, 9).Value = 997.7143
$ws.Cells.Item(22, 11).Value = 997.7143
$ws.Cells.Item(22, 13).Value = -647.7143

$ws.Cells.Item(31, 8).Value = 2225.6191
$ws.Cells.Item(31, 9).Value = 2018
$ws.Cells.Item(31, 11).Value = 2018
$ws.Cells.Item(31, 13).Value = -1723

$ws.Cells.Item(34, 8).Value = 2225.6191
$ws.Cells.Item(34, 9).Value = 2018
$ws.Cells.Item(34, 11).Value = 2018
$ws.Cells.Item(34, 13).Value = -1816

$ws.Cells.Item(41, 9).Value = 0
$ws.Cells.Item(41, 10).Value = 20000
$ws.Cells.Item(41, 11).Value = 0
$ws.Cells.Item(41, 12).Value = 20000
$ws.Cells.Item(41, 14).Value = -20856
$ws.Cells.Item(41, 13).ClearContents()

$ws.Cells.Item(58, 8).Value = 1033.5
$ws.Cells.Item(58, 9).Value = 1133.4166
$ws.Cells.Item(58, 10).Value = 733.75
$ws.Cells.Item(58, 11).Value = 1133.4166
$ws.Cells.Item(58, 12).Value = 733.75
$ws.Cells.Item(58, 13).Value = -930.4166
$ws.Cells.Item(58, 14).Value = -1139.75

$ws.Cells.Item(86, 8).Value = 6226.25
$ws.Cells.Item(86, 10).Value = 4954
$ws.Cells.Item(86, 12).Value = 4954
$ws.Cells.Item(86, 14).Value = -7200

$ws.Cells.Item(87, 8).Value = 61666.668
$ws.Cells.Item(87, 10).Value = 61666.668
$ws.Cells.Item(87, 12).Value = 61666.668
$ws.Cells.Item(87, 14).Value = -64038.668

$ws.Cells.Item(89, 8).Value = 6226.25
$ws.Cells.Item(89, 10).Value = 4954
$ws.Cells.Item(89, 12).Value = 24770
$ws.Cells.Item(89, 14).Value = -36002

$ws.Cells.Item(90, 8).Value = 61666.668
$ws.Cells.Item(90, 10).Value = 61666.668
$ws.Cells.Item(90, 12).Value = 185000.004
$ws.Cells.Item(90, 14).Value = -196856.004

$ws.Cells.Item(94, 8).Value = 674.7778
$ws.Cells.Item(94, 9).Value = 807
$ws.Cells.Item(94, 11).Value = 807
$ws.Cells.Item(94, 13).Value = -356

$ws.Cells.Item(99, 8).Value = 8003.5
$ws.Cells.Item(99, 10).Value = 9337.666999999999
$ws.Cells.Item(99, 12).Value = 9337.666999999999
$ws.Cells.Item(99, 14).Value = -12333.667

$ws.Cells.Item(105, 8).Value = 794
$ws.Cells.Item(105, 9).Value = 794
$ws.Cells.Item(105, 11).Value = 794
$ws.Cells.Item(105, 13).Value = 953

$ws.Cells.Item(107, 8).Value = 1900.2572
$ws.Cells.Item(107, 9).Value = 1354.6957
$ws.Cells.Item(107, 10).Value = 2945.9167
$ws.Cells.Item(107, 11).Value = 1354.6957
$ws.Cells.Item(107, 12).Value = 2945.9167
$ws.Cells.Item(107, 13).Value = 565.3043
$ws.Cells.Item(107, 14).Value = -6785.9167

$ws.Cells.Item(124, 8).Value = 89996.336
$ws.Cells.Item(124, 10).Value = 89996.336
$ws.Cells.Item(124, 12).Value = 89996.336
$ws.Cells.Item(124, 14).Value = -94906.336

$ws.Cells.Item(126, 8).Value = 8003.5
$ws.Cells.Item(126, 10).Value = 9337.666999999999
$ws.Cells.Item(126, 12).Value = 28013.001
$ws.Cells.Item(126, 14).Value = -32953.001

$ws.Cells.Item(132, 8).Value = 2716.1785
$ws.Cells.Item(132, 9).Value = 2574.3914
$ws.Cells.Item(132, 10).Value = 3368.4
$ws.Cells.Item(132, 11).Value = 7723.174199999999
$ws.Cells.Item(132, 12).Value = 10105.2
$ws.Cells.Item(132, 13).Value = -5193.174199999999
$ws.Cells.Item(132, 14).Value = -15165.2

$ws.Cells.Item(134, 8).Value = 2086.5
$ws.Cells.Item(134, 9).Value = 1691.8478
$ws.Cells.Item(134, 11).Value = 5075.5434
$ws.Cells.Item(134, 13).Value = -2540.5434

$ws.Cells.Item(136, 8).Value = 1033.5
$ws.Cells.Item(136, 9).Value = 1133.4166
$ws.Cells.Item(136, 10).Value = 733.75
$ws.Cells.Item(136, 11).Value = 3400.2498
$ws.Cells.Item(136, 12).Value = 2201.25
$ws.Cells.Item(136, 13).Value = -850.2498000000001
$ws.Cells.Item(136, 14).Value = -7301.25

$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(5, 8).Value = 52498.5
$ws.Cells.Item(5, 9).Value = 0
$ws.Cells.Item(5, 10).Value = 52498.5
$ws.Cells.Item(5, 11).Value = 0
$ws.Cells.Item(5, 12).Value = 157495.5
$ws.Cells.Item(5, 14).Value = -157719.5
$ws.Cells.Item(5, 13).ClearContents()

$ws.Cells.Item(17, 8).Value = 86.833336
$ws.Cells.Item(17, 9).Value = 78.666664
$ws.Cells.Item(17, 10).Value = 95
$ws.Cells.Item(17, 11).Value = 235.999992
$ws.Cells.Item(17, 12).Value = 285
$ws.Cells.Item(17, 13).Value = -66.99999199999999
$ws.Cells.Item(17, 14).Value = -623

$ws.Cells.Item(18, 8).Value = 3155.2222
$ws.Cells.Item(18, 9).Value = 3124.5
$ws.Cells.Item(18, 10).Value = 3179.8
$ws.Cells.Item(18, 11).Value = 9373.5
$ws.Cells.Item(18, 12).Value = 9539.400000000001
$ws.Cells.Item(18, 13).Value = -9204.5
$ws.Cells.Item(18, 14).Value = -9877.400000000001

$ws.Cells.Item(25, 8).Value = 1219.625
$ws.Cells.Item(25, 9).Value = 328.57144
$ws.Cells.Item(25, 10).Value = 1912.6666
$ws.Cells.Item(25, 11).Value = 985.71432
$ws.Cells.Item(25, 12).Value = 5737.9998
$ws.Cells.Item(25, 13).Value = -816.71432
$ws.Cells.Item(25, 14).Value = -6075.9998

$ws.Cells.Item(29, 8).Value = 2554.9
$ws.Cells.Item(29, 9).Value = 2499.5
$ws.Cells.Item(29, 10).Value = 2568.75
$ws.Cells.Item(29, 11).Value = 7498.5
$ws.Cells.Item(29, 12).Value = 7706.25
$ws.Cells.Item(29, 14).Value = -8260.25
$ws.Cells.Item(29, 13).Value = -7221.5

$ws.Cells.Item(30, 8).Value = 1219.625
$ws.Cells.Item(30, 9).Value = 328.57144
$ws.Cells.Item(30, 10).Value = 1912.6666
$ws.Cells.Item(30, 11).Value = 985.71432
$ws.Cells.Item(30, 12).Value = 5737.9998
$ws.Cells.Item(30, 13).Value = -883.71432
$ws.Cells.Item(30, 14).Value = -5941.9998

$ws.Cells.Item(36, 8).Value = 182.8
$ws.Cells.Item(36, 9).Value = 182.8
$ws.Cells.Item(36, 11).Value = 548.4000000000001
$ws.Cells.Item(36, 13).Value = -379.4000000000001

$ws.Cells.Item(47, 8).Value = 775
$ws.Cells.Item(47, 10).Value = 833.3333
$ws.Cells.Item(47, 12).Value = 2499.9999
$ws.Cells.Item(47, 14).Value = -3361.9999

$ws.Cells.Item(56, 8).Value = 2387533.8
$ws.Cells.Item(56, 9).Value = 2387533.8
$ws.Cells.Item(56, 11).Value = 2387533.8
$ws.Cells.Item(56, 13).Value = -2387003.8

$ws.Cells.Item(98, 8).Value = 856.2143
$ws.Cells.Item(98, 9).Value = 906.55554
$ws.Cells.Item(98, 11).Value = 2719.66662
$ws.Cells.Item(98, 13).Value = -1221.66662

$ws.Cells.Item(114, 8).Value = 3807.5
$ws.Cells.Item(114, 10).Value = 5810.2
$ws.Cells.Item(114, 12).Value = 17430.6
$ws.Cells.Item(114, 14).Value = -23938.6

$ws.Cells.Item(129, 8).Value = 1966.0476
$ws.Cells.Item(129, 9).Value = 925.63635
$ws.Cells.Item(129, 10).Value = 3110.5
$ws.Cells.Item(129, 11).Value = 2776.90905
$ws.Cells.Item(129, 12).Value = 9331.5
$ws.Cells.Item(129, 13).Value = 2223.09095
$ws.Cells.Item(129, 14).Value = -19331.5

$ws.Cells.Item(131, 8).Value = 1645.5084
$ws.Cells.Item(131, 10).Value = 1712.7037
$ws.Cells.Item(131, 12).Value = 5138.1111
$ws.Cells.Item(131, 14).Value = -15218.1111

$ws.Cells.Item(135, 8).Value = 52498.5
$ws.Cells.Item(135, 9).Value = 0
$ws.Cells.Item(135, 10).Value = 52498.5
$ws.Cells.Item(135, 11).Value = 0
$ws.Cells.Item(135, 12).Value = 472486.5
$ws.Cells.Item(135, 14).Value = -477556.5
$ws.Cells.Item(135, 13).ClearContents()

$ws.Cells.Item(137, 8).Value = 6183.1333
$ws.Cells.Item(137, 9).Value = 4224.6665
$ws.Cells.Item(137, 10).Value = 7488.778
$ws.Cells.Item(137, 11).Value = 12673.9995
$ws.Cells.Item(137, 12).Value = 22466.334
$ws.Cells.Item(137, 13).Value = -7573.999500000002
$ws.Cells.Item(137, 14).Value = -32666.334

$ws.Cells.Item(140, 8).Value = 3066.2
$ws.Cells.Item(140, 9).Value = 1978.3
$ws.Cells.Item(140, 10).Value = 5242
$ws.Cells.Item(140, 11).Value = 5934.9
$ws.Cells.Item(140, 12).Value = 15726
$ws.Cells.Item(140, 13).Value = -754.8999999999996
$ws.Cells.Item(140, 14).Value = -26086

$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(3, 8).Value = 11879750
$ws.Cells.Item(3, 10).Value = 15837000
$ws.Cells.Item(3, 12).Value = 15837000
$ws.Cells.Item(3, 14).Value = -15837232

$ws.Cells.Item(19, 8).Value = 19666.666
$ws.Cells.Item(19, 10).Value = 0
$ws.Cells.Item(19, 12).Value = 0
$ws.Cells.Item(19, 14).ClearContents()

$ws.Cells.Item(43, 8).Value = 61395.312
$ws.Cells.Item(43, 9).Value = 22056
$ws.Cells.Item(43, 10).Value = 84998.89999999999
$ws.Cells.Item(43, 11).Value = 22056
$ws.Cells.Item(43, 12).Value = 84998.89999999999
$ws.Cells.Item(43, 13).Value = -21905
$ws.Cells.Item(43, 14).Value = -85300.89999999999

$ws.Cells.Item(52, 8).Value = 23499.6
$ws.Cells.Item(52, 10).Value = 23749.5
$ws.Cells.Item(52, 12).Value = 23749.5
$ws.Cells.Item(52, 14).Value = -24267.5

$ws.Cells.Item(70, 8).Value = 3996.6155
$ws.Cells.Item(70, 9).Value = 2795.2
$ws.Cells.Item(70, 11).Value = 2795.2
$ws.Cells.Item(70, 13).Value = -2525.2

$ws.Cells.Item(73, 8).Value = 3996.6155
$ws.Cells.Item(73, 9).Value = 2795.2
$ws.Cells.Item(73, 11).Value = 2795.2
$ws.Cells.Item(73, 13).Value = -1859.2

$ws.Cells.Item(97, 8).Value = 1181.08
$ws.Cells.Item(97, 9).Value = 903.7857
$ws.Cells.Item(97, 10).Value = 1534
$ws.Cells.Item(97, 11).Value = 903.7857
$ws.Cells.Item(97, 12).Value = 1534
$ws.Cells.Item(97, 13).Value = -407.7857
$ws.Cells.Item(97, 14).Value = -2526

$ws.Cells.Item(101, 8).Value = 60461.832
$ws.Cells.Item(101, 10).Value = 60461.832
$ws.Cells.Item(101, 12).Value = 60461.832
$ws.Cells.Item(101, 14).Value = -66951.83199999999

$ws.Cells.Item(102, 8).Value = 43589.645
$ws.Cells.Item(102, 9).Value = 3475.7
$ws.Cells.Item(102, 10).Value = 143874.5
$ws.Cells.Item(102, 11).Value = 3475.7
$ws.Cells.Item(102, 12).Value = 143874.5
$ws.Cells.Item(102, 13).Value = -1853.7
$ws.Cells.Item(102, 14).Value = -147118.5

$ws.Cells.Item(122, 8).Value = 1385.4242
$ws.Cells.Item(122, 9).Value = 1229.6774
$ws.Cells.Item(122, 10).Value = 3799.5
$ws.Cells.Item(122, 11).Value = 3689.0322
$ws.Cells.Item(122, 12).Value = 11398.5
$ws.Cells.Item(122, 13).Value = -1239.0322
$ws.Cells.Item(122, 14).Value = -16298.5

$ws.Cells.Item(123, 8).Value = 34199.2
$ws.Cells.Item(123, 10).Value = 34199.2
$ws.Cells.Item(123, 12).Value = 34199.2
$ws.Cells.Item(123, 14).Value = -39099.2

$ws.Cells.Item(126, 8).Value = 2635
$ws.Cells.Item(126, 9).Value = 2427.818
$ws.Cells.Item(126, 11).Value = 7283.454000000001
$ws.Cells.Item(126, 13).Value = -4813.454000000001

$ws.Cells.Item(132, 8).Value = 7514.1904
$ws.Cells.Item(132, 9).Value = 8824.875
$ws.Cells.Item(132, 10).Value = 3320
$ws.Cells.Item(132, 11).Value = 26474.625
$ws.Cells.Item(132, 12).Value = 9960
$ws.Cells.Item(132, 13).Value = -23944.625
$ws.Cells.Item(132, 14).Value = -15020

$ws.Cells.Item(139, 8).Value = 149999.5
$ws.Cells.Item(139, 10).Value = 149999.5
$ws.Cells.Item(139, 12).Value = 149999.5
$ws.Cells.Item(139, 14).Value = -160279.5

$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(7, 8).Value = 3234.318
$ws.Cells.Item(7, 9).Value = 1814.3889
$ws.Cells.Item(7, 10).Value = 9624
$ws.Cells.Item(7, 11).Value = 1814.3889
$ws.Cells.Item(7, 12).Value = 9624
$ws.Cells.Item(7, 13).Value = -1702.3889
$ws.Cells.Item(7, 14).Value = -9848

$ws.Cells.Item(22, 8).Value = 12785.077
$ws.Cells.Item(22, 10).Value = 27226.666
$ws.Cells.Item(22, 12).Value = 27226.666
$ws.Cells.Item(22, 14).Value = -27816.666

$ws.Cells.Item(26, 8).Value = 15000
$ws.Cells.Item(26, 10).Value = 15000
$ws.Cells.Item(26, 12).Value = 15000
$ws.Cells.Item(26, 14).Value = -15590

$ws.Cells.Item(27, 8).Value = 12785.077
$ws.Cells.Item(27, 10).Value = 27226.666
$ws.Cells.Item(27, 12).Value = 27226.666
$ws.Cells.Item(27, 14).Value = -27440.666

$ws.Cells.Item(30, 8).Value = 7850
$ws.Cells.Item(30, 9).Value = 7850
$ws.Cells.Item(30, 10).Value = 0
$ws.Cells.Item(30, 11).Value = 7850
$ws.Cells.Item(30, 12).Value = 0
$ws.Cells.Item(30, 13).Value = -7742
$ws.Cells.Item(30, 14).ClearContents()

$ws.Cells.Item(31, 8).Value = 9000
$ws.Cells.Item(31, 9).Value = 3000
$ws.Cells.Item(31, 10).Value = 15000
$ws.Cells.Item(31, 11).Value = 3000
$ws.Cells.Item(31, 12).Value = 15000
$ws.Cells.Item(31, 13).Value = -2752
$ws.Cells.Item(31, 14).Value = -15496

$ws.Cells.Item(40, 8).Value = 2602.7666
$ws.Cells.Item(40, 9).Value = 2331.6072
$ws.Cells.Item(40, 11).Value = 2331.6072
$ws.Cells.Item(40, 13).Value = -2195.6072

$ws.Cells.Item(46, 8).Value = 8582.777
$ws.Cells.Item(46, 10).Value = 14450
$ws.Cells.Item(46, 12).Value = 14450
$ws.Cells.Item(46, 14).Value = -14826

$ws.Cells.Item(61, 8).Value = 1355.9286
$ws.Cells.Item(61, 9).Value = 1280.5
$ws.Cells.Item(61, 11).Value = 1280.5
$ws.Cells.Item(61, 13).Value = -1078.5

$ws.Cells.Item(95, 8).Value = 75000
$ws.Cells.Item(95, 10).Value = 75000
$ws.Cells.Item(95, 12).Value = 75000
$ws.Cells.Item(95, 14).Value = -80492

$ws.Cells.Item(99, 8).Value = 23000
$ws.Cells.Item(99, 9).Value = 23000
$ws.Cells.Item(99, 11).Value = 23000
$ws.Cells.Item(99, 13).Value = -20005

$ws.Cells.Item(113, 8).Value = 1355.9286
$ws.Cells.Item(113, 9).Value = 1280.5
$ws.Cells.Item(113, 11).Value = 1280.5
$ws.Cells.Item(113, 13).Value = 889.5

$ws.Cells.Item(122, 8).Value = 8386.546
$ws.Cells.Item(122, 9).Value = 7419.3335
$ws.Cells.Item(122, 11).Value = 22258.0005
$ws.Cells.Item(122, 13).Value = -19808.0005

$ws.Cells.Item(126, 8).Value = 3234.318
$ws.Cells.Item(126, 9).Value = 1814.3889
$ws.Cells.Item(126, 10).Value = 9624
$ws.Cells.Item(126, 11).Value = 5443.1667
$ws.Cells.Item(126, 12).Value = 28872
$ws.Cells.Item(126, 13).Value = -2973.1667
$ws.Cells.Item(126, 14).Value = -33812

$ws.Cells.Item(132, 8).Value = 1873.7778
$ws.Cells.Item(132, 9).Value = 1385.6774
$ws.Cells.Item(132, 10).Value = 4900
$ws.Cells.Item(132, 11).Value = 4157.0322
$ws.Cells.Item(132, 12).Value = 14700
$ws.Cells.Item(132, 13).Value = -1627.0322
$ws.Cells.Item(132, 14).Value = -19760

$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(15, 8).Value = 0
$ws.Cells.Item(15, 10).Value = 0
$ws.Cells.Item(15, 12).Value = 0
$ws.Cells.Item(15, 14).ClearContents()

$ws.Cells.Item(53, 8).Value = 84
$ws.Cells.Item(53, 10).Value = 84
$ws.Cells.Item(53, 12).Value = 84
$ws.Cells.Item(53, 14).Value = -1298

$ws.Cells.Item(56, 8).Value = 43438
$ws.Cells.Item(56, 10).Value = 50157
$ws.Cells.Item(56, 12).Value = 50157
$ws.Cells.Item(56, 14).Value = -51585

$ws.Cells.Item(58, 8).Value = 40833.332
$ws.Cells.Item(58, 10).Value = 40833.332
$ws.Cells.Item(58, 12).Value = 40833.332
$ws.Cells.Item(58, 14).Value = -41449.332

$ws.Cells.Item(62, 8).Value = 16800586
$ws.Cells.Item(62, 9).Value = 199630.75
$ws.Cells.Item(62, 11).Value = 199630.75
$ws.Cells.Item(62, 13).Value = -199006.75

$ws.Cells.Item(65, 8).Value = 16800586
$ws.Cells.Item(65, 9).Value = 199630.75
$ws.Cells.Item(65, 11).Value = 998153.75
$ws.Cells.Item(65, 13).Value = -995033.75

$ws.Cells.Item(95, 8).Value = 74999
$ws.Cells.Item(95, 10).Value = 74999
$ws.Cells.Item(95, 12).Value = 74999
$ws.Cells.Item(95, 14).Value = -80491

$ws.Cells.Item(100, 8).Value = 1424.5
$ws.Cells.Item(100, 9).Value = 1344.5238
$ws.Cells.Item(100, 10).Value = 1760.4
$ws.Cells.Item(100, 11).Value = 2689.0476
$ws.Cells.Item(100, 12).Value = 3520.8
$ws.Cells.Item(100, 13).Value = -2148.0476
$ws.Cells.Item(100, 14).Value = -4602.8

$ws.Cells.Item(107, 8).Value = 324.54285
$ws.Cells.Item(107, 9).Value = 246.15
$ws.Cells.Item(107, 11).Value = 738.45
$ws.Cells.Item(107, 13).Value = 1181.55

$ws.Cells.Item(122, 8).Value = 2674.2415
$ws.Cells.Item(122, 9).Value = 2409.3704
$ws.Cells.Item(122, 11).Value = 7228.111199999999
$ws.Cells.Item(122, 13).Value = -4778.111199999999

$ws.Cells.Item(126, 8).Value = 1830.3636
$ws.Cells.Item(126, 9).Value = 1652.3
$ws.Cells.Item(126, 10).Value = 3611
$ws.Cells.Item(126, 11).Value = 4956.9
$ws.Cells.Item(126, 12).Value = 10833
$ws.Cells.Item(126, 13).Value = -2486.9
$ws.Cells.Item(126, 14).Value = -15773

$ws.Cells.Item(132, 8).Value = 1298.1406
$ws.Cells.Item(132, 9).Value = 1286.0807
$ws.Cells.Item(132, 10).Value = 1672
$ws.Cells.Item(132, 11).Value = 3858.2421
$ws.Cells.Item(132, 12).Value = 5016
$ws.Cells.Item(132, 13).Value = -1328.2421
$ws.Cells.Item(132, 14).Value = -10076
